$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: Column Subset changed from "all" to "subrun overlap fts"
$ws.Range("C29").Value = "subrun overlap fts"

# Row 34: fill in K..O (job results for an already-described run)
$ws.Range("K34").Value = 161
$ws.Range("L34").Value = "34.5 & 21.4"
$ws.Range("M34").Value = "58.8 & 43.7"
$ws.Range("N34").Value = 15
$ws.Range("O34").Value = 35.9

# Row 35: new job row
$ws.Range("A35").Value = "ukb51139_subset.csv"
$ws.Range("B35").Value = "2801 x 147"
$ws.Range("C35").Value = "subrun overlap fts"
$ws.Range("D35").Value = "no events"
$ws.Range("E35").Value = "> 140/80"
$ws.Range("F35").Value = "zscore"
$ws.Range("G35").Value = "median"
$ws.Range("H35").Value = "none"
$ws.Range("I35").Value = 50
$ws.Range("K35").Value = 49
$ws.Range("L35").Value = "81.6 & 77.4"
$ws.Range("M35").Value = "68.3 & 63.7"
$ws.Range("N35").Value = 14
$ws.Range("O35").Value = 2.65

# Row 36: new job row
$ws.Range("A36").Value = "ukb51139_subset.csv"
$ws.Range("B36").Value = "7003 x 147"
$ws.Range("C36").Value = "subrun overlap fts"
$ws.Range("D36").Value = "no events"
$ws.Range("E36").Value = "> 140/80"
$ws.Range("F36").Value = "zscore"
$ws.Range("G36").Value = "median"
$ws.Range("H36").Value = "none"
$ws.Range("I36").Value = 50
$ws.Range("K36").Value = 46
$ws.Range("L36").Value = "95.8 & 94.2"
$ws.Range("M36").Value = "79.9 & 75.2"
$ws.Range("N36").Value = 18
$ws.Range("O36").Value = 3.2

# Row 37: new job row
$ws.Range("A37").Value = "ukb51139_subset.csv"
$ws.Range("B37").Value = "4902 x 462"
$ws.Range("C37").Value = "subrun overlap fts"
$ws.Range("D37").Value = "no events"
$ws.Range("E37").Value = "> 140/80"
$ws.Range("F37").Value = "zscore"
$ws.Range("G37").Value = "median"
$ws.Range("H37").Value = "none"
$ws.Range("I37").Value = 50
$ws.Range("K37").Value = 63
$ws.Range("L37").Value = "55.0 & 40.7"
$ws.Range("M37").Value = "48.9 & 43.1"
$ws.Range("N37").Value = 18
$ws.Range("O37").Value = 32.7

# Row 38: new job row
$ws.Range("A38").Value = "ukb51139_subset.csv"
$ws.Range("B38").Value = "3502 x 1081"
$ws.Range("C38").Value = "all"
$ws.Range("D38").Value = "no events"
$ws.Range("E38").Value = "> 140/80"
$ws.Range("F38").Value = "zscore"
$ws.Range("G38").Value = "median"
$ws.Range("H38").Value = "none"
$ws.Range("I38").Value = 50
$ws.Range("K38").Value = 480
$ws.Range("L38").Value = "-398.3 & -40.7"
$ws.Range("M38").Value = "35.3 & 34.9"
$ws.Range("N38").Value = 16
$ws.Range("O38").Value = 68.9
